# Applies the "Updated symbol list" data refresh to the crypto tracking sheet.
# All data cells in this sheet are stored as text (inline strings), including
# cells whose text happens to look like a number (e.g. "284.62"). When such a
# value is written through COM with a plain .Value assignment, Excel silently
# re-interprets it as a numeric cell and normalizes its representation (e.g.
# dropping meaningful trailing zeros such as "0.03220" -> 0.0322). To avoid
# that, numeric-looking values are written while the cell is temporarily
# forced to Text format, which keeps Excel from re-typing the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.NumberFormat = "General"
}

function Set-PlainCell {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $ws.Range($Address).Value = $NewValue
}

# --- Price (column D) updates throughout the coin list ---
Set-TextCell "D2"  "284.62"
Set-TextCell "D4"  "6.449"
Set-TextCell "D5"  "0.06349"
Set-TextCell "D6"  "3.598"
Set-TextCell "D7"  "1.531"
Set-TextCell "D8"  "6.559"
Set-TextCell "D9"  "0.8204"
Set-TextCell "D10" "0.01412"
Set-TextCell "D11" "0.1673"
Set-TextCell "D12" "0.08621"
Set-TextCell "D13" "0.03657"
Set-TextCell "D14" "0.03220"
Set-TextCell "D15" "0.09198"
Set-TextCell "D16" "3.723"
Set-TextCell "D17" "0.001648"
Set-TextCell "D18" "0.04745"
Set-TextCell "D19" "0.006152"
Set-TextCell "D20" "0.006274"
Set-TextCell "D24" "2.270"
Set-TextCell "D25" "0.3355"
Set-TextCell "D40" "0.04762"

# --- Rows 41-43: coin ranking list rotated (Kick/BKEX/CEJI reshuffled) ---
Set-PlainCell "B41" "KickToken"
Set-PlainCell "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell  "D41" "0.007103"
Set-PlainCell "E41" "40KickTokenKICK"

Set-PlainCell "B42" "BKEXToken"
Set-PlainCell "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell  "D42" "0.1114"
Set-PlainCell "E42" "41BKEXTokenBKK"

Set-PlainCell "B43" "CEJI"
Set-PlainCell "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell  "D43" "0.003967"
Set-PlainCell "E43" "42CEJICEJI"

# --- Remaining price updates ---
Set-TextCell "D44" "0.01144"
Set-TextCell "D45" "0.00007100"
Set-TextCell "D47" "1.002"
Set-TextCell "D48" "0.003988"
Set-TextCell "D49" "0.00001503"
Set-PlainCell "E49" "48CryptobidCoinCBCWorstin24h"
